$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.429.37"
$ws.Range("E2").Value = "  +0.66%  "
$ws.Range("D3").Value = "1.612.37"
$ws.Range("E3").Value = "  +1.23%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.88"
$ws.Range("E5").Value = "  -0.71%  "
$ws.Range("E6").Value = "  -0.36%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -0.41%  "
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("E10").Value = "  +1.39%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0846"
$ws.Range("E11").Value = "  -0.59%  "
$ws.Range("E12").Value = "  +1.28%  "
$ws.Range("D13").Value = "1.609.41"
$ws.Range("E13").Value = "  +1.18%  "
$ws.Range("E14").Value = "  -0.19%  "
$ws.Range("E15").Value = "  +0.05%  "
$ws.Range("E16").Value = "  -0.28%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "234.27"
$ws.Range("E17").Value = "  +8.80%  "
$ws.Range("D18").Value = "26.429.25"
$ws.Range("E18").Value = "  +0.69%  "
$ws.Range("D19").Value = "0.0₃0725"
$ws.Range("E19").Value = "  +0.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.66"
$ws.Range("E20").Value = "  +3.97%  "
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("E22").Value = "  -0.27%  "
$ws.Range("E23").Value = "  +4.43%  "
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.69"
$ws.Range("E25").Value = "  +1.32%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("E27").Value = "  +0.28%  "
$ws.Range("E28").Value = "  +0.26%  "
$ws.Range("E29").Value = "  +2.30%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0496"
$ws.Range("E30").Value = "  +1.12%  "
$ws.Range("E31").Value = "  -0.48%  "
$ws.Range("D32").Value = "1.497.35"
$ws.Range("E32").Value = "  +6.16%  "
$ws.Range("E33").Value = "  +1.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.94"
$ws.Range("E34").Value = "  -1.03%  "
$ws.Range("E36").Value = "  +2.40%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.560"
$ws.Range("E37").Value = "  -2.56%  "
$ws.Range("E38").Value = "  -0.06%  "
$ws.Range("E39").Value = "  +0.42%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.80"
$ws.Range("E40").Value = "  +0.64%  "
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.18"
$ws.Range("E42").Value = "  +1.06%  "
$ws.Range("D43").Value = "1.751.50"
$ws.Range("E43").Value = "  +1.37%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.926"
$ws.Range("E44").Value = "  -3.29%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.761"
$ws.Range("E45").Value = "  -0.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "61.31"
$ws.Range("E46").Value = "  +0.72%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "89.64"
$ws.Range("E47").Value = "  +2.83%  "
$ws.Range("E48").Value = "  -0.25%  "
$ws.Range("E49").Value = "  -0.11%  "
$ws.Range("E50").Value = "  +0.94%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.47"
$ws.Range("E51").Value = "  +1.01%  "
